# Finish chain from csv to a hashtable of perform categories
#
# The sheet holds "skill categories": a table that maps a trigger
# (coords/type/scope/direction) to an impact (coords/type/scope/direction).
# Previously there was no "impact coords" column, so the impact side of the
# table was missing a step that the trigger side already had. This adds the
# missing "impact coords" column (mirroring "trigger coords") and re-derives
# the impact type/scope values for the "attack"/"arrow" rows, which used to
# be flagged "single" (no scope) and now use "scope" with an explicit value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook was imported from a CSV, so reflect that in the sheet title.
$ws.Name = "skill_categories.csv"

# Insert a new column G ("impact coords"), shifting the old impact
# type/scope/direction columns (G,H,I) right to (H,I,J) - mirrors the
# trigger coords/type/scope/direction layout already present in C..F.
$ws.Columns("G").Insert()
$ws.Range("G1").Value = "impact coords"
$ws.Columns("G").ColumnWidth = 12

# Rows 2 & 3 ("attack", "arrow") previously had impact type "single" with no
# impact scope recorded. Finish the chain: they resolve to impact type
# "scope" with an explicit scope value of 1 (a single target).
$ws.Range("H2").Value = "scope"
$ws.Range("I2").Value = 1

$ws.Range("H3").Value = "scope"
$ws.Range("I3").Value = 1

# Rows 4 & 5 ("shoot", "canon") already carried correct impact type/scope
# pairs - the column insert above shifted them from G/H into H/I intact, so
# nothing further to change there.

# Leave the selection where the last edit landed, matching the authored file.
$ws.Range("H5").Select() | Out-Null
